# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets,
# matching the new generated snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F3"  = 12877
        "F5"  = 86
        "F6"  = 70
        "F9"  = 12
        "F10" = 12798
        "F13" = 8667
        "F15" = 190
        "F16" = 94
        "F20" = 7
        "F25" = 85
    }
    "全部类型" = @{
        "F4"  = 12877
        "F6"  = 86
        "F7"  = 70
        "F10" = 12
        "F11" = 12798
        "F14" = 8667
        "F16" = 190
        "F17" = 94
        "F21" = 7
        "F27" = 85
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
